$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column widths (A..AA) ---
$widths = @(10,8,15,11,32,18,4,20,23,24,20,7,17,11,26,26,27,32,13,31,27,21,33,31,42,20,28)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i] - 0.8333333333333334
}

# --- Freeze header row (pane split after row 1) ---
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$null = $ws.Range("A1").Select()

# --- Header row formatting: center/center + wrap text on A1:Z1 ---
$hdr = $ws.Range("A1:Z1")
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true

# --- AA1: same alignment + yellow fill + updated text ---
$aa1 = $ws.Range("AA1")
$aa1.HorizontalAlignment = -4108
$aa1.VerticalAlignment = -4108
$aa1.WrapText = $true
$aa1.Interior.Color = 65535
$aa1.Value = "Status as of July 11, 2025"

# --- AA2:AA7: yellow fill on the dropdown status column ---
$aaBody = $ws.Range("AA2:AA7")
$aaBody.Interior.Color = 65535
